# RASPBERRY-3B-it_eval.xlsx update
#
# Task #12 ("extract 5 topics") was re-scored: the VOTE cell had been a
# stray text value "3t" (so it didn't count in the SUBTOTAL) and the
# REMARKS cell said "here are only 3". Fix the vote to the real numeric
# score (3) and correct the remark's wording to "there are only 3".
# The TOTAL (SUBTOTAL formula) recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F18 = VOTE for row #12 ("extract 5 topics"): was the text "3t", should
# be the number 3 so it is included in the SUBTOTAL total below.
$ws.Range("F18").Value = 3

# G18 = REMARKS for the same row: correct the typo/wording.
$ws.Range("G18").Value = "there are only 3"

# The last touched/selected cell in the sheet moved to A12.
$ws.Range("A12").Select()
